# Update the division answers in the practice table to the new output
# generated at 9a8706d.
#
# Every cell's old text is unique within the document, so a sequence of
# exact (non-wildcard) Find/Replace operations over the whole document
# content unambiguously targets the correct cell for each substitution.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "96÷6=16, 0" "38÷5=7, 3"
Replace-Text "53÷7=7, 4" "39÷2=19, 1"
Replace-Text "30÷5=6, 0" "79÷6=13, 1"
Replace-Text "24÷6=4, 0" "98÷2=49, 0"
Replace-Text "46÷2=23, 0" "91÷5=18, 1"
Replace-Text "94÷7=13, 3" "15÷6=2, 3"
Replace-Text "64÷4=16, 0" "87÷3=29, 0"
Replace-Text "54÷7=7, 5" "39÷7=5, 4"
Replace-Text "60÷3=20, 0" "14÷4=3, 2"
Replace-Text "33÷7=4, 5" "23÷7=3, 2"
Replace-Text "37÷8=4, 5" "33÷9=3, 6"
Replace-Text "45÷8=5, 5" "50÷9=5, 5"
Replace-Text "42÷8=5, 2" "46÷5=9, 1"
# Note: "94÷3=31, 1" is both the existing 5th cell of this row AND the
# value that the (soon to be replaced) "68÷3=22, 2" cell will become, so
# this replacement must run before that one to avoid double-matching.
Replace-Text "94÷3=31, 1" "80÷9=8, 8"
Replace-Text "68÷3=22, 2" "94÷3=31, 1"
Replace-Text "53÷6=8, 5" "37÷7=5, 2"
Replace-Text "54÷8=6, 6" "86÷7=12, 2"
Replace-Text "94÷6=15, 4" "69÷6=11, 3"
Replace-Text "55÷8=6, 7" "72÷9=8, 0"
Replace-Text "42÷7=6, 0" "99÷6=16, 3"
Replace-Text "57÷4=14, 1" "95÷9=10, 5"
Replace-Text "89÷9=9, 8" "80÷5=16, 0"
Replace-Text "17÷9=1, 8" "85÷6=14, 1"
Replace-Text "53÷4=13, 1" "19÷8=2, 3"
Replace-Text "45÷4=11, 1" "87÷5=17, 2"
